# Inserts a new weekly price record for Coliflor (Vega Modelo de Temuco)
# as row 406, pushing the existing rows 406-437 down to 407-438.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 406. Excel will shift all
# rows 406..437 down to 407..438 and carry the date-cell style (column D)
# down from the row above, matching the original layout.
$ws.Rows.Item(406).Insert()

# Populate the newly inserted row 406 with the new record.
$ws.Cells.Item(406, 1).Value = 10
$ws.Cells.Item(406, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(406, 3).Value = "La Araucanía"
$ws.Cells.Item(406, 4).Value = 44783
$ws.Cells.Item(406, 5).Value = 9
$ws.Cells.Item(406, 6).Value = 100112008
$ws.Cells.Item(406, 7).Value = "Coliflor"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Primera"
$ws.Cells.Item(406, 10).Value = 800
$ws.Cells.Item(406, 11).Value = 1300
$ws.Cells.Item(406, 12).Value = 1300
$ws.Cells.Item(406, 13).Value = 1300
$ws.Cells.Item(406, 14).Value = "`$/unidad"
$ws.Cells.Item(406, 15).Value = "Región Metropolitana"
$ws.Cells.Item(406, 16).Value = 1300
$ws.Cells.Item(406, 17).Value = 1
$ws.Cells.Item(406, 18).Value = "Hortaliza"
